$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style used by the
# other header cells (bold, bordered, centered) - copy format from H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for the new columns I and J, rows 2-22
$data = @(
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(14, 14),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(7, 8),
    @(8, 9),
    @(1, 4),
    @(6, 7),
    @(8, 8),
    @(4, 6),
    @(4, 4),
    @(6, 6),
    @(5, 5),
    @(3, 3),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
